# Generate Report for Handback
# Updates handoff/handback timestamps for the d4e93497 file across the
# zh-cn and de-de status sheets, and refreshes the "Latest HO Xliff
# Generate Date" summary column on the Overview sheet.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsOverview = $wb.Worksheets.Item("Overview")

# zh-cn: row 3 corresponds to d4e93497-4b60-497b-a193-6f002878760b.md
$wsZhCn.Range("H3").Value = "2016-08-29 22:50:52"
$wsZhCn.Range("K3").Value = "2016-08-29 22:51:18"

# de-de: row 3 corresponds to d4e93497-4b60-497b-a193-6f002878760b.md
$wsDeDe.Range("H3").Value = "2016-08-29 22:50:56"
$wsDeDe.Range("K3").Value = "2016-08-29 22:51:25"

# Overview: row 3 corresponds to d4e93497-4b60-497b-a193-6f002878760b.md,
# column G holds the latest handoff xliff generation date across languages.
$wsOverview.Range("G3").Value = "2016-08-29 22:50:56"
